# Insert a new worksheet "Solvay Plant" right after "brine", summarising the
# process steps of the Solvay plant, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$brine = $wb.Worksheets.Item("brine")
$newSheet = $wb.Worksheets.Add($null, $brine)
$newSheet.Name = "Solvay Plant"

$newSheet.Range("A2").Value = "ammonia absorber"
$newSheet.Range("A3").Value = "filter 1"
$newSheet.Range("A4").Value = "solvay tower"
$newSheet.Range("A5").Value = "lime kiln"
$newSheet.Range("A6").Value = "slaker"

# Match the author's selection/cursor position and make this the active tab.
$newSheet.Activate()
$newSheet.Range("A7").Select()
